# Port_Letter.xlsx — "mergeFromTo array settings added, some minor refactor added"
#
# The sheet gets a brand new, almost-invisible spacer row inserted right
# above the signature block (old row 41), pushing everything from the old
# row 41 down through row 52 one row further down. The previously-blank
# F-cell of the (now shifted) comment row gets a "-" placeholder, a couple
# of named ranges get re-pointed to the shifted cells, and two new named
# ranges (Merge_end / Pg_end) are introduced to support the new
# mergeFromTo array logic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert the new spacer row at row 41 (shifts old rows 41-52 -> 42-53)
# ---------------------------------------------------------------------
$ws.Rows("41:41").Insert()

# New row 41 is a near-invisible divider holding a single "-" in A41.
$ws.Range("A41").Value = "-"
$ws.Rows("41:41").RowHeight = 0.85

# ---------------------------------------------------------------------
# 2. Old row 43 (now row 44) gains a "-" placeholder in its F cell,
#    which used to be blank.
# ---------------------------------------------------------------------
$ws.Range("F44").Value = "-"

# ---------------------------------------------------------------------
# 3. The conditional formatting range that used to end exactly on the
#    inserted row (H29:N41) needs to grow along with it -> H29:N42.
# ---------------------------------------------------------------------
$cf = $ws.Cells.FormatConditions.Item(2)
$cf.ModifyAppliesToRange($ws.Range("H29:N42"))

# ---------------------------------------------------------------------
# 4. Re-point the named ranges that lived below row 40 to their new,
#    shifted-down locations.
# ---------------------------------------------------------------------
$wb.Names.Item("Seal_seller_start").RefersTo = "=Port_Letter!`$D`$42"
$wb.Names.Item("Seal_seller_end").RefersTo = "=Port_Letter!`$E`$44"
$wb.Names.Item("Sign_seller_start").RefersTo = "=Port_Letter!`$D`$43"
$wb.Names.Item("Sign_seller_end").RefersTo = "=Port_Letter!`$E`$43"
$wb.Names.Item("Подписант").RefersTo = "=Port_Letter!`$F`$43"
$wb.Names.Item("Подписант_комментарий").RefersTo = "=Port_Letter!`$A`$43"

# New named ranges introduced for the mergeFromTo array settings.
$wb.Names.Add("Merge_end", "=Port_Letter!`$A`$41")
$wb.Names.Add("Pg_end", "=Port_Letter!`$F`$44")

# ---------------------------------------------------------------------
# 5. Print area grows by the extra row.
# ---------------------------------------------------------------------
$ws.PageSetup.PrintArea = "`$A`$1:`$F`$52"

# ---------------------------------------------------------------------
# 6. Restore the author's on-screen selection.
# ---------------------------------------------------------------------
$ws.Range("C40").Select() | Out-Null
